# Update the "想去人数" (F column) counts for both the "展览" and "全部类型"
# sheets to reflect freshly generated output (gh-pages rebuild at 456a3b4).

$wb = $excel.ActiveWorkbook

# Row number (as shown in column A / sheet row) -> new F-column value
$updates = @{
    2  = 11729
    3  = 11357
    11 = 10771
    12 = 4165
    18 = 51
    20 = 451
    21 = 11144
    22 = 10929
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
